$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: split single "שם" column into first/last name headers
$ws.Range("A1").Value = "שם פרטי"
$ws.Range("B1").Value = "שם משפחה"

# Data rows: split "last first..." full name into first-name (col A) and
# last-name (col B) parts. The rows below contain, in order:
#   column A = first name (can be multiple words)
#   column B = last name  (can be multiple words)
$names = @(
    @("ריף", "אדרי"),
    @("ליאור", "בנדרסקי"),
    @("בניה חיים", "גרובר"),
    @("יותם", "הנדורגר"),
    @("יהלי", "זפרני"),
    @("עדן", "טיומקין"),
    @("הראל אנריקה", "טייב"),
    @("יוני", "ילין לנדסקרו"),
    @("מאיה", "ישראל"),
    @("אורי דוד", "כחלון"),
    @("דור", "לנדמן"),
    @("מאי", "סלע"),
    @("הראל", "פסטמן"),
    @("עמית", "רובין"),
    @("דנאל", "שוסטרמן"),
    @("יונתן", "שיינברג"),
    @("עמית", "שטופמכר")
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i][0]
    $ws.Cells.Item($row, 2).Value = $names[$i][1]
}

# Set column widths for the new layout (both columns A and B).
# Target OOXML width is 15.83203125 characters; the COM ColumnWidth setter
# snaps to whole-pixel increments, so 15 (which rounds to the same pixel
# bucket) is the closest achievable value.
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 15
